$d = $word.ActiveDocument
$quoteOpen  = [char]8220
$quoteClose = [char]8221

# ---------------------------------------------------------------------------
# Helper: collapse an entire paragraph's runs into a single run without
# changing its visible text (forces the engine to re-merge adjoining runs
# and drop spell-check proofErr markers, mirroring what Word does when the
# text in a paragraph is reselected/retyped).
# ---------------------------------------------------------------------------
function Merge-Paragraph($para) {
    $rng = $para.Range.Duplicate()
    $rng.MoveEnd(1, -1)   # exclude the paragraph mark
    $orig = $rng.Text
    $rng.Text = "ZZZ_MERGE_PLACEHOLDER_ZZZ"
    $rng2 = $para.Range.Duplicate()
    $rng2.MoveEnd(1, -1)
    $rng2.Text = $orig
}

# ---------------------------------------------------------------------------
# Helper: re-establish a run boundary immediately before/around a specific
# piece of text inside an already-merged paragraph, by toggling a
# formatting property on and back off (a no-op visually, but it forces the
# engine to split that span back out into its own run).
# $searchText  - the literal text of the span that should become its own run
# $context     - extra text immediately following $searchText used only to
#                disambiguate the Find match (not included in the split)
# ---------------------------------------------------------------------------
function Split-Off($para, $searchText, $context) {
    $rng = $para.Range.Duplicate()
    $rng.Find.ClearFormatting()
    $full = $searchText + $context
    $rng.Find.Execute($full, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $rng.End = $rng.Start + $searchText.Length
    $rng.Bold = 1
    $rng.Bold = 0
}

$paras = $d.Paragraphs

# --- Paragraph 4: "In script / rmd file "01-TrackingData-Introduction-Visualise"
#     , need to properly source and cite Yelkouan Shearwater tracking data
#     example from BIOM" -------------------------------------------------
$p = $paras.Item(4)
Merge-Paragraph $p
Split-Off $p "01-TrackingData-Introduction-Visualise" ""
Split-Off $p $quoteClose ", need"

# --- Paragraph 6: "In script / rmd file "01-TrackingData-Introduction-Visualise""
$p = $paras.Item(6)
Merge-Paragraph $p
Split-Off $p "01-TrackingData-Introduction-Visualise" ""

# --- Paragraph 11: "In script / rmd file "01-TrackingData-Introduction-Visualise""
$p = $paras.Item(11)
Merge-Paragraph $p
Split-Off $p "01-TrackingData-Introduction-Visualise" ""

# --- Paragraph 16: "In script / rmd "02-TrackingData-SamplingStrategy""
$p = $paras.Item(16)
Merge-Paragraph $p
Split-Off $p "02-TrackingData-SamplingStrategy" ""

# --- Paragraph 19: "In script / rmd "02-TrackingData-SamplingStrategy""
$p = $paras.Item(19)
Merge-Paragraph $p
Split-Off $p "02-TrackingData-SamplingStrategy" ""

# --- Paragraph 22: "In script / rmd "03-TrackingData-Visualisation"" (full collapse)
$p = $paras.Item(22)
Merge-Paragraph $p

# --- Paragraph 23: "Need to create figure showcasing Lastovo SPA and source
#     population for Yelkouan Shearwaters" (keep trailing "s" as its own run)
$p = $paras.Item(23)
Merge-Paragraph $p
$rngS = $p.Range.Duplicate()
$rngS.MoveEnd(1, -1)
$rngS.Start = $rngS.End - 1
$rngS.Bold = 1
$rngS.Bold = 0

# --- Paragraph 24: "In script / rmd "03-TrackingData-Visualisation"" (full collapse)
$p = $paras.Item(24)
Merge-Paragraph $p

# --- Paragraph 28: "In script / rmd "03-TrackingData-Visualisation"" (full collapse)
$p = $paras.Item(28)
Merge-Paragraph $p

# --- Paragraph 30: "In script / rmd "03-TrackingData-Visualisation"" (full collapse)
$p = $paras.Item(30)
Merge-Paragraph $p

# ---------------------------------------------------------------------------
# Append the three new bullet paragraphs after the
# "[Decide on best way to show example datasets ...]" paragraph.
# ---------------------------------------------------------------------------
$lastPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Decide on best way to show example datasets*") {
        $lastPara = $para
    }
}

$rngEnd = $lastPara.Range.Duplicate()
$rngEnd.Collapse(0)
$rngEnd.InsertParagraphAfter()

# Re-fetch paragraphs collection / locate the three freshly inserted (still
# empty) paragraphs that follow $lastPara.
$newPara1 = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Start -eq $rngEnd.Start) {
        $newPara1 = $para
    }
}
# Fallback: locate by scanning since Start offsets shift after insertion.
$idx = 0
$targetIdx = -1
foreach ($para in $d.Paragraphs) {
    $idx = $idx + 1
    if ($para.Range.Text -like "*Decide on best way to show example datasets*") {
        $targetIdx = $idx
    }
}
$newPara1 = $d.Paragraphs.Item($targetIdx + 1)
$newPara1.Range.ListFormat.ListLevelNumber = 1
$newPara1.Range.Text = "In script / rmd " + $quoteOpen + "03-TrackingData-Visualisation" + $quoteClose

$rngEnd2 = $newPara1.Range.Duplicate()
$rngEnd2.Collapse(0)
$rngEnd2.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item($targetIdx + 2)
$newPara2.Range.ListFormat.ListLevelNumber = 2
$newPara2.Range.Text = "> Visualising tracking data from individual animals can help you understand which data you might remove, or which data you might try and salvage."

$rngEnd3 = $newPara2.Range.Duplicate()
$rngEnd3.Collapse(0)
$rngEnd3.InsertParagraphAfter()
$newPara3 = $d.Paragraphs.Item($targetIdx + 3)
$newPara3.Range.ListFormat.ListLevelNumber = 3
$newPara3.Range.Text = "Need to provide examples of tracks you would remove, and which tracks you might try and salvage"

Write-Host "Done"
